$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename product description "ROMA MARI GOLD110+22" -> "ROMA MARI GOLD 110G"
$ws.Range("B186").Value = "ROMA MARI GOLD 110G"

# Remove the discontinued product row (20078242 / GERY MLK SLT CNUT105),
# shifting all subsequent rows up by one.
$ws.Rows(211).Delete()
